$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 239. This shifts the existing rows
# 239-277 down to 240-278, preserving all their data and formatting
# (including the date number format applied to column D).
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with the new data record.
$ws.Cells.Item(239, 1).Value = 10
$ws.Cells.Item(239, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(239, 3).Value = "La Araucanía"
$ws.Cells.Item(239, 4).Value = 44491
$ws.Cells.Item(239, 5).Value = 9
$ws.Cells.Item(239, 6).Value = 100112023
$ws.Cells.Item(239, 7).Value = "Brócoli"
$ws.Cells.Item(239, 8).Value = "Sin especificar"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 2770
$ws.Cells.Item(239, 11).Value = 800
$ws.Cells.Item(239, 12).Value = 900
$ws.Cells.Item(239, 13).Value = 844
$ws.Cells.Item(239, 14).Value = "`$/unidad"
$ws.Cells.Item(239, 15).Value = "Región Metropolitana"
$ws.Cells.Item(239, 16).Value = 844
$ws.Cells.Item(239, 17).Value = 1
$ws.Cells.Item(239, 18).Value = "Hortaliza"
